$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 156 (pushes the existing rows 156:188 down to 157:189)
$ws.Rows("156:156").Insert()

# Populate the newly inserted row with the new weekly price-report entry
$ws.Range("A156").Value2 = 10
$ws.Range("B156").Value2 = "Vega Modelo de Temuco"
$ws.Range("C156").Value2 = "La Araucanía"
$ws.Range("D156").Value2 = 44637
$ws.Range("E156").Value2 = 9
$ws.Range("F156").Value2 = 100112005
$ws.Range("G156").Value2 = "Puerro"
$ws.Range("H156").Value2 = "Azul de Maquehue"
$ws.Range("I156").Value2 = "Primera"
$ws.Range("J156").Value2 = 40
$ws.Range("K156").Value2 = 10000
$ws.Range("L156").Value2 = 10000
$ws.Range("M156").Value2 = 10000
$ws.Range("N156").Value2 = "$/docena de paquetes"
$ws.Range("O156").Value2 = "Provincia de Cautín"
$ws.Range("P156").Value2 = 833
$ws.Range("Q156").Value2 = 12
$ws.Range("R156").Value2 = "Hortaliza"
